# Week 16 logging + season sim update
$wb = $excel.ActiveWorkbook
$wsR = $wb.Worksheets.Item("Rushing")
$wsC = $wb.Worksheets.Item("Receiving")

# ---------------------------------------------------------------------------
# Rushing sheet: update existing rows, then insert the new player N.Cottrell
# into row 7 (pushing the old row7-9 data that follows down), and append a
# new T.Austin rushing row at the end.
# ---------------------------------------------------------------------------

# Row 2 - T.Lawrence
$wsR.Range("C2").Value = 22
$wsR.Range("D2").Value = 12
$wsR.Range("E2").Value = 23
$wsR.Range("F2").Value = 12

# Row 4 - J.Robinson
$wsR.Range("C4").Value = 92
$wsR.Range("E4").Value = 11
$wsR.Range("F4").Value = 24

# Row 6 - D.Ogunbowale
$wsR.Range("C6").Value = 12
$wsR.Range("D6").Value = 8
$wsR.Range("E6").Value = 3
$wsR.Range("F6").Value = 3

# Row 7 - new player N.Cottrell (was L.Shenault's slot; L.Shenault's own row
# moves to row 8 below)
$wsR.Range("B7").Value = "N.Cottrell"
$wsR.Range("C7").Value = 0
$wsR.Range("D7").Value = 1
$wsR.Range("E7").Value = 0
$wsR.Range("F7").Value = 0

# Row 8 - L.Shenault
$wsR.Range("B8").Value = "L.Shenault"
$wsR.Range("C8").Value = 4
$wsR.Range("D8").Value = 7
$wsR.Range("E8").Value = 0
$wsR.Range("F8").Value = 2

# Row 9 - T.Johnson
$wsR.Range("B9").Value = "T.Johnson"
$wsR.Range("C9").Value = 1
$wsR.Range("D9").Value = 0
$wsR.Range("E9").Value = 0
$wsR.Range("F9").Value = 1

# Row 10 (new) - J.Agnew, formatted like the existing index/name rows
$wsR.Range("A9").Copy()
$wsR.Range("A10").PasteSpecial(-4122)
$wsR.Range("A10").Value = 8
$wsR.Range("B10").Value = "J.Agnew"
$wsR.Range("C10").Value = 5
$wsR.Range("D10").Value = 2
$wsR.Range("E10").Value = 0
$wsR.Range("F10").Value = 0

# Row 11 (new) - T.Austin
$wsR.Range("A9").Copy()
$wsR.Range("A11").PasteSpecial(-4122)
$wsR.Range("A11").Value = 9
$wsR.Range("B11").Value = "T.Austin"
$wsR.Range("C11").Value = 2
$wsR.Range("D11").Value = 1
$wsR.Range("E11").Value = 0
$wsR.Range("F11").Value = 1

# ---------------------------------------------------------------------------
# Receiving sheet: roster names are unaffected, only season totals change.
# ---------------------------------------------------------------------------

# Row 4 - D.Ogunbowale
$wsC.Range("C4").Value = 13
$wsC.Range("D4").Value = 8
$wsC.Range("G4").Value = 2

# Row 5 - M.Jones
$wsC.Range("C5").Value = 75
$wsC.Range("D5").Value = 54
$wsC.Range("E5").Value = 30
$wsC.Range("G5").Value = 13
$wsC.Range("H5").Value = 5

# Row 9 - T.Austin
$wsC.Range("C9").Value = 28
$wsC.Range("D9").Value = 19
$wsC.Range("E9").Value = 6
$wsC.Range("F9").Value = 2

# Row 10 - L.Treadwell
$wsC.Range("C10").Value = 25
$wsC.Range("D10").Value = 19
$wsC.Range("E10").Value = 12
$wsC.Range("F10").Value = 5

# Row 13 - J.O'Shaughnessy
$wsC.Range("C13").Value = 27
$wsC.Range("D13").Value = 19
$wsC.Range("E13").Value = 8
$wsC.Range("F13").Value = 5

# Row 14 - J.Hollister
$wsC.Range("C14").Value = 9
$wsC.Range("D14").Value = 6

Write-Host "Week 16 stats applied"
